# Append 4 more trial/iteration rows to the "Data" sheet (rows 6-9) and the
# corresponding compact rows to the "Summary" sheet (rows 5-8), for a new
# SHAPESHIFTER run whose chromosome string is recorded once and shared.

$wb = $excel.ActiveWorkbook
$wsData    = $wb.Worksheets.Item("Data")
$wsSummary = $wb.Worksheets.Item("Summary")

# New chromosome / individual description shared by all 4 new Data rows.
$chromosome = "[9, 10, 15, 14, 128, 0.07882309, 0.23540762, 0.6077191, 0.87509125, 0.57467437, 0.17408052, 0.8004514, 0.6979225, 0.033451915, 0.24232545, 0.68438506, 0.8240594, 0.4639135, 0.5343747, 0.7130215, 0.31194702, 0.50155693, 0.49896312, 0.28083193, 0.4807646, 0.37663388, 0.9757447, 0.28397393, 0.986866, 0.553169, 0.28753248, 0.7259604, 0.19067651, 0.90213567, 0.6042983, 0.39408544, 0.27314886, 0.85486674, 0.6418463, 0.37925383, 0.3937522, 0.78596526, 0.5726512, 0.42413598, 0.67200243, 0.032528486, 0.134553, 0.48343706, 0.5637112, 0.73201555, 0.66249174, 0.269676]"

# Columns B, E-G, H (chromosome), I-AD stay constant for every new Data row;
# only A (iteration), C and D (fitness / std.dev. that keep shrinking) vary.
$constant = @{
     2 = 2.0333333;    5 = 2.0333333;    6 = 6.7337713;   7 = 30.0;
     9 = 361.7;        10 = 135.66621;   11 = 358.23334;  12 = 19.786243;
    13 = 130.53334;    14 = 53.297302;   15 = 0.033333335; 16 = 0.18257418;
    17 = 30.9;         18 = 7.438715;    19 = 30.9;        20 = 7.438715;
    21 = 16.066668;    22 = 7.3152337;   23 = 0.0;         24 = 0.0;
    25 = 30.9;         26 = 7.438715;    27 = 2.0333333;   28 = 6.7337713;
    29 = 114.166664;   30 = 26.07692
}

# Per-row varying values: iteration (A), fitness (C), std.dev. (D)
$newRows = @(
    @{ Row = 6; A = 3.0; C = 0.067777775;    D = 0.37123418   },
    @{ Row = 7; A = 4.0; C = 0.0033333334;   D = 0.018257419  },
    @{ Row = 8; A = 5.0; C = 0.0022222223;   D = 0.012171612  },
    @{ Row = 9; A = 6.0; C = 0.0;            D = 0.0          }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $wsData.Cells.Item($row, 1).Value = $r.A
    $wsData.Cells.Item($row, 3).Value = $r.C
    $wsData.Cells.Item($row, 4).Value = $r.D
    $wsData.Cells.Item($row, 8).Value = $chromosome

    foreach ($col in $constant.Keys) {
        $wsData.Cells.Item($row, $col).Value = $constant[$col]
    }

    # Mirror the same A/B/C (iteration, avg time, fitness) onto the Summary
    # sheet, four rows up (Data row 6 -> Summary row 5, etc.).
    $summaryRow = $row - 1
    $wsSummary.Cells.Item($summaryRow, 1).Value = $r.A
    $wsSummary.Cells.Item($summaryRow, 2).Value = $constant[2]
    $wsSummary.Cells.Item($summaryRow, 3).Value = $r.C
}
